# Apply cryptocurrency price/volume updates scraped on Thu Jul 27 10:34:19 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.555.26"
$ws.Range("E2").Value = "  +1.23%  "
$ws.Range("D3").Value = "1.882.34"
$ws.Range("E3").Value = "  +1.61%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'0.7163"
$ws.Range("E5").Value = "  +2.56%  "
$ws.Range("D6").Value = "'242.30"
$ws.Range("E6").Value = "  +1.98%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "'0.07947"
$ws.Range("D9").Value = "'0.3116"
$ws.Range("E9").Value = "  +3.36%  "
$ws.Range("E10").Value = "  +6.99%  "
$ws.Range("D11").Value = "'0.08283"
$ws.Range("E11").Value = "  +2.19%  "
$ws.Range("D12").Value = "'0.7310"
$ws.Range("E12").Value = "  +3.83%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'5.292"
$ws.Range("E13").Value = "  +2.26%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.862.47"
$ws.Range("E14").Value = "  +0.46%  "
$ws.Range("D15").Value = "'91.42"
$ws.Range("E15").Value = "  +2.19%  "
$ws.Range("D16").Value = "29.535.63"
$ws.Range("E16").Value = "  +1.03%  "
$ws.Range("D17").Value = "'5.949"
$ws.Range("E17").Value = "  +2.60%  "
$ws.Range("D18").Value = "'246.41"
$ws.Range("E18").Value = "  +4.59%  "
$ws.Range("D19").Value = "'0.000007896"
$ws.Range("E19").Value = "  +1.04%  "
$ws.Range("D20").Value = "'13.37"
$ws.Range("E20").Value = "  +1.41%  "
$ws.Range("D21").Value = "2.133.43"
$ws.Range("E21").Value = "  +1.37%  "
$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").Value = "'7.979"
$ws.Range("E23").Value = "  +6.47%  "
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").Value = "'0.1618"
$ws.Range("E25").Value = "  +13.74%  "
$ws.Range("D26").Value = "'163.40"
$ws.Range("E26").Value = "  +0.58%  "
$ws.Range("D27").Value = "'9.074"
$ws.Range("E27").Value = "  +2.46%  "
$ws.Range("D28").Value = "'18.35"
$ws.Range("E28").Value = "  +1.90%  "
$ws.Range("D29").Value = "'1.359"
$ws.Range("E29").Value = "  -3.26%  "
$ws.Range("D30").Value = "'1.502"
$ws.Range("E30").Value = "  +1.75%  "
$ws.Range("D31").Value = "'4.398"
$ws.Range("E31").Value = "  +1.72%  "
$ws.Range("D32").Value = "'4.115"
$ws.Range("E32").Value = "  +2.64%  "
$ws.Range("D33").Value = "'0.05307"
$ws.Range("E33").Value = "  +2.90%  "
$ws.Range("D34").Value = "'1.964"
$ws.Range("D35").Value = "'1.203"
$ws.Range("E35").Value = "  +3.48%  "
$ws.Range("D36").Value = "'0.7282"
$ws.Range("E36").Value = "  +2.49%  "
$ws.Range("D37").Value = "'2.681"
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("D38").Value = "'0.01873"
$ws.Range("E38").Value = "  +1.52%  "
$ws.Range("D39").Value = "1.233.43"
$ws.Range("E39").Value = "  +7.30%  "
$ws.Range("E40").Value = "  +0.33%  "
$ws.Range("D41").Value = "'0.9117"
$ws.Range("E41").Value = "  -1.16%  "
$ws.Range("D42").Value = "'6.215"
$ws.Range("E42").Value = "  +4.05%  "
$ws.Range("D43").Value = "'73.88"
$ws.Range("E43").Value = "  +5.51%  "
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("D45").Value = "'102.22"
$ws.Range("E45").Value = "  -0.68%  "
$ws.Range("D46").Value = "2.028.99"
$ws.Range("E46").Value = "  +1.38%  "
$ws.Range("D47").Value = "'0.5285"
$ws.Range("E47").Value = "  -0.11%  "
$ws.Range("D48").Value = "'1.798"
$ws.Range("E48").Value = "  +3.67%  "
$ws.Range("D49").Value = "'2.945"
$ws.Range("E49").Value = "  +11.50%  "
$ws.Range("E50").Value = "  +2.62%  "
$ws.Range("D51").Value = "'9.330"
$ws.Range("E51").Value = "  +2.22%  "
